$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(8, 1)
$run = $para.Runs(1)
$run.Text = "Ensemble methods (Random Forests, Boosting, Bagging, etc)"
